$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new "remaining days" E value, new "start date" F value)
# Daily countdown update: each row's remaining-days counter (column E) ticks down by
# one. When a counter would hit zero, it resets to the row's total-days value (column D)
# and the start date (column F) rolls forward by that same number of days.
$updates = @{
    2 = @(2, 20260126)
    3 = @(2, 20260126)
    4 = @(2, 20260126)
    5 = @(6, 20260203)
    6 = @(2, 20260126)
    7 = @(6, 20260203)
    8 = @(2, 20260126)
    9 = @(6, 20260203)
    10 = @(2, 20260202)
    11 = @(2, 20260126)
    12 = @(6, 20260203)
    13 = @(2, 20260126)
    14 = @(2, 20260126)
    15 = @(2, 20260126)
    16 = @(10, 20260207)
    17 = @(6, 20260203)
    18 = @(9, 20260206)
    19 = @(9, 20260206)
    20 = @(9, 20260206)
    21 = @(9, 20260206)
    22 = @(6, 20260203)
    23 = @(6, 20260203)
    24 = @(6, 20260203)
    25 = @(6, 20260203)
    26 = @(6, 20260203)
    27 = @(3, 20260203)
    28 = @(9, 20260206)
    29 = @(9, 20260206)
    30 = @(9, 20260206)
    31 = @(9, 20260206)
    32 = @(9, 20260206)
    33 = @(9, 20260206)
    34 = @(9, 20260206)
    35 = @(9, 20260206)
    37 = @(9, 20260206)
    38 = @(9, 20260206)
    39 = @(9, 20260206)
    40 = @(2, 20260202)
    41 = @(2, 20260202)
    42 = @(9, 20260206)
    43 = @(6, 20260203)
    44 = @(2, 20260202)
    45 = @(6, 20260203)
    46 = @(2, 20260202)
    47 = @(9, 20260206)
    48 = @(2, 20260202)
    49 = @(3, 20260203)
    50 = @(4, 20260201)
    51 = @(4, 20260201)
    52 = @(4, 20260201)
    53 = @(4, 20260201)
    54 = @(4, 20260201)
    55 = @(4, 20260201)
    56 = @(4, 20260201)
    57 = @(4, 20260201)
    58 = @(8, 20260205)
    59 = @(8, 20260205)
    60 = @(8, 20260205)
    61 = @(3, 20260203)
    62 = @(8, 20260205)
    63 = @(8, 20260205)
    64 = @(8, 20260205)
    65 = @(9, 20260206)
    66 = @(9, 20260206)
    67 = @(9, 20260206)
    68 = @(9, 20260206)
    69 = @(9, 20260206)
    70 = @(10, 20260207)
    71 = @(10, 20260207)
    72 = @(10, 20260207)
    73 = @(10, 20260207)
    74 = @(10, 20260207)
    75 = @(10, 20260207)
    76 = @(10, 20260207)
    77 = @(3, 20260131)
    78 = @(3, 20260131)
    79 = @(3, 20260131)
    80 = @(3, 20260131)
    81 = @(3, 20260131)
    82 = @(3, 20260131)
    83 = @(3, 20260131)
    84 = @(3, 20260131)
    85 = @(3, 20260131)
    86 = @(3, 20260131)
    87 = @(2, 20260202)
    88 = @(2, 20260202)
    89 = @(2, 20260202)
    90 = @(2, 20260202)
    91 = @(6, 20260203)
    92 = @(2, 20260202)
    93 = @(3, 20260131)
    94 = @(5, 20260205)
    95 = @(2, 20260130)
    96 = @(10, 20260207)
    97 = @(10, 20260207)
    98 = @(10, 20260207)
    99 = @(10, 20260207)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]   # column E - 剩余 (remaining days)
    $ws.Cells.Item($row, 6).Value = $vals[1]   # column F - 开始时间 (start date, yyyymmdd)
}
